$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font, border, centered alignment) of the existing
# header cell H1 onto the two new header cells I1 and J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

# New data values for row 3
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 8
